$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: Overview
#   Existing row 3 (b68a2137) must become row 5.
#   New rows 3 and 4 are inserted for 10bb3719 and 66575674.
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Preserve old row3 values before overwriting.
$old3_A = $ws1.Range("A3").Value2
$old3_B = $ws1.Range("B3").Value2
$old3_C = $ws1.Range("C3").Value2
$old3_D = $ws1.Range("D3").Value2

# Drop existing hyperlinks (cell text/values are untouched); they will be rebuilt below.
$ws1.Hyperlinks.Delete()

# Row 5 <- old row 3 data (b68a2137)
$ws1.Range("A5").Value2 = $old3_A
$ws1.Range("B5").Value2 = $old3_B
$ws1.Range("C5").Value2 = $old3_C
$ws1.Range("D5").Value2 = $old3_D
$ws1.Range("D5").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Row 3 <- 10bb3719
$ws1.Range("A3").Value2 = "10bb3719-a9cc-4ae1-86d9-c76b788262a8.md"
$ws1.Range("B3").Value2 = "Ready for handoff"
$ws1.Range("C3").Value2 = "Ready for handoff"
$ws1.Range("D3").Value2 = "2016-03-24 04:40:08"
$ws1.Range("D3").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Row 4 <- 66575674
$ws1.Range("A4").Value2 = "66575674-fd33-4ba1-9e47-8dae58941a75.md"
$ws1.Range("B4").Value2 = "Ready for handoff"
$ws1.Range("C4").Value2 = "Ready for handoff"
$ws1.Range("D4").Value2 = "2016-03-24 04:40:08"
$ws1.Range("D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

# Rebuild hyperlinks for A2..A5
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/6e684c28b587a70f0cfd5f5c6ff1ff3adb27a702/e2e/20bcb66c-dd1a-43b3-a248-76b294441b45.md", "", "", "20bcb66c-dd1a-43b3-a248-76b294441b45.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/10bb3719a9cc4ae186d9c76b788262a8000000000/e2e/10bb3719-a9cc-4ae1-86d9-c76b788262a8.md", "", "", "10bb3719-a9cc-4ae1-86d9-c76b788262a8.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/66575674fd334ba19e478dae58941a75000000000/e2e/66575674-fd33-4ba1-9e47-8dae58941a75.md", "", "", "66575674-fd33-4ba1-9e47-8dae58941a75.md")
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/4a68d5de10ad6c3300541df7be93d85829774398/e2e/b68a2137-e2a6-439b-88bb-be13c4dbe045.md", "", "", "b68a2137-e2a6-439b-88bb-be13c4dbe045.md")

# ---------------------------------------------------------------
# Sheet 2: zh-cn
#   Same row-shift pattern as sheet 1, but with more columns
#   (A,B,C,D,E,H,J) and two hyperlinked columns (A and D).
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$old3_A = $ws2.Range("A3").Value2
$old3_B = $ws2.Range("B3").Value2
$old3_C = $ws2.Range("C3").Value2
$old3_D = $ws2.Range("D3").Value2
$old3_E = $ws2.Range("E3").Value2
$old3_H = $ws2.Range("H3").Value2
$old3_J = $ws2.Range("J3").Value2

$ws2.Hyperlinks.Delete()

# Row 5 <- old row 3 data (b68a2137)
$ws2.Range("A5").Value2 = $old3_A
$ws2.Range("B5").Value2 = $old3_B
$ws2.Range("C5").Value2 = $old3_C
$ws2.Range("D5").Value2 = $old3_D
$ws2.Range("E5").Value2 = $old3_E
$ws2.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H5").Value2 = $old3_H
$ws2.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("J5").Value2 = $old3_J

# Row 3 <- 10bb3719
$ws2.Range("A3").Value2 = "10bb3719-a9cc-4ae1-86d9-c76b788262a8.md"
$ws2.Range("B3").Value2 = ".md"
$ws2.Range("C3").Value2 = "Ready for handoff"
$ws2.Range("D3").Value2 = "10bb3719-a9cc-4ae1-86d9-c76b788262a8.b35df4560331639d247d8c3db9bcf659d9491bc4.zh-cn.xlf"
$ws2.Range("E3").Value2 = "2016-03-24 04:40:02"
$ws2.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H3").Value2 = "0001-01-01 00:00:00"
$ws2.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("J3").Value2 = "Include"

# Row 4 <- 66575674
$ws2.Range("A4").Value2 = "66575674-fd33-4ba1-9e47-8dae58941a75.md"
$ws2.Range("B4").Value2 = ".md"
$ws2.Range("C4").Value2 = "Ready for handoff"
$ws2.Range("D4").Value2 = "66575674-fd33-4ba1-9e47-8dae58941a75.2d11a3a5a73b03f5b090dc2164ab42aedd19db64.zh-cn.xlf"
$ws2.Range("E4").Value2 = "2016-03-24 04:40:02"
$ws2.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H4").Value2 = "0001-01-01 00:00:00"
$ws2.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("J4").Value2 = "Include"

# Rebuild hyperlinks for A2,D2,F2,G2 (unchanged) then A3,D3,A4,D4,A5,D5
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/6e684c28b587a70f0cfd5f5c6ff1ff3adb27a702/e2e/20bcb66c-dd1a-43b3-a248-76b294441b45.md", "", "", "20bcb66c-dd1a-43b3-a248-76b294441b45.md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/579fedc10ff59a7149311f4f12d601c1d77743ec/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/20bcb66c-dd1a-43b3-a248-76b294441b45.633f0c36d9a87b0eefdf682d49f1e52b4ef346d6.zh-cn.xlf", "", "", "20bcb66c-dd1a-43b3-a248-76b294441b45.633f0c36d9a87b0eefdf682d49f1e52b4ef346d6.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/9a05043f1ca8cbd194aa7a8cceb5a0d1fef77f3a/e2e/20bcb66c-dd1a-43b3-a248-76b294441b45.md", "", "", "20bcb66c-dd1a-43b3-a248-76b294441b45.md")
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/76b17313277728859431b60dbd28f02cc35e369d/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/20bcb66c-dd1a-43b3-a248-76b294441b45.633f0c36d9a87b0eefdf682d49f1e52b4ef346d6.zh-cn.xlf", "", "", "20bcb66c-dd1a-43b3-a248-76b294441b45.633f0c36d9a87b0eefdf682d49f1e52b4ef346d6.zh-cn.xlf")

$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/10bb3719a9cc4ae186d9c76b788262a8000000000/e2e/10bb3719-a9cc-4ae1-86d9-c76b788262a8.md", "", "", "10bb3719-a9cc-4ae1-86d9-c76b788262a8.md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b35df4560331639d247d8c3db9bcf659d9491bc4/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/10bb3719-a9cc-4ae1-86d9-c76b788262a8.b35df4560331639d247d8c3db9bcf659d9491bc4.zh-cn.xlf", "", "", "10bb3719-a9cc-4ae1-86d9-c76b788262a8.b35df4560331639d247d8c3db9bcf659d9491bc4.zh-cn.xlf")

$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/66575674fd334ba19e478dae58941a75000000000/e2e/66575674-fd33-4ba1-9e47-8dae58941a75.md", "", "", "66575674-fd33-4ba1-9e47-8dae58941a75.md")
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2d11a3a5a73b03f5b090dc2164ab42aedd19db64/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/66575674-fd33-4ba1-9e47-8dae58941a75.2d11a3a5a73b03f5b090dc2164ab42aedd19db64.zh-cn.xlf", "", "", "66575674-fd33-4ba1-9e47-8dae58941a75.2d11a3a5a73b03f5b090dc2164ab42aedd19db64.zh-cn.xlf")

$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/4a68d5de10ad6c3300541df7be93d85829774398/e2e/b68a2137-e2a6-439b-88bb-be13c4dbe045.md", "", "", "b68a2137-e2a6-439b-88bb-be13c4dbe045.md")
$ws2.Hyperlinks.Add($ws2.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/7cf28389d94558a735df52069639c0c7a33fa611/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b68a2137-e2a6-439b-88bb-be13c4dbe045.cafe1cadf14af19aedb092ec9298155126eb95ed.zh-cn.xlf", "", "", "b68a2137-e2a6-439b-88bb-be13c4dbe045.cafe1cadf14af19aedb092ec9298155126eb95ed.zh-cn.xlf")

# ---------------------------------------------------------------
# Sheet 3: de-de
#   Same row-shift pattern, de-de variant.
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$old3_A = $ws3.Range("A3").Value2
$old3_B = $ws3.Range("B3").Value2
$old3_C = $ws3.Range("C3").Value2
$old3_D = $ws3.Range("D3").Value2
$old3_E = $ws3.Range("E3").Value2
$old3_H = $ws3.Range("H3").Value2
$old3_J = $ws3.Range("J3").Value2

$ws3.Hyperlinks.Delete()

# Row 5 <- old row 3 data (b68a2137)
$ws3.Range("A5").Value2 = $old3_A
$ws3.Range("B5").Value2 = $old3_B
$ws3.Range("C5").Value2 = $old3_C
$ws3.Range("D5").Value2 = $old3_D
$ws3.Range("E5").Value2 = $old3_E
$ws3.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H5").Value2 = $old3_H
$ws3.Range("H5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("J5").Value2 = $old3_J

# Row 3 <- 10bb3719
$ws3.Range("A3").Value2 = "10bb3719-a9cc-4ae1-86d9-c76b788262a8.md"
$ws3.Range("B3").Value2 = ".md"
$ws3.Range("C3").Value2 = "Ready for handoff"
$ws3.Range("D3").Value2 = "10bb3719-a9cc-4ae1-86d9-c76b788262a8.b35df4560331639d247d8c3db9bcf659d9491bc4.de-de.xlf"
$ws3.Range("E3").Value2 = "2016-03-24 04:40:08"
$ws3.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H3").Value2 = "0001-01-01 00:00:00"
$ws3.Range("H3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("J3").Value2 = "Include"

# Row 4 <- 66575674
$ws3.Range("A4").Value2 = "66575674-fd33-4ba1-9e47-8dae58941a75.md"
$ws3.Range("B4").Value2 = ".md"
$ws3.Range("C4").Value2 = "Ready for handoff"
$ws3.Range("D4").Value2 = "66575674-fd33-4ba1-9e47-8dae58941a75.2d11a3a5a73b03f5b090dc2164ab42aedd19db64.de-de.xlf"
$ws3.Range("E4").Value2 = "2016-03-24 04:40:08"
$ws3.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H4").Value2 = "0001-01-01 00:00:00"
$ws3.Range("H4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("J4").Value2 = "Include"

# Rebuild hyperlinks for A2,D2,F2,G2 (unchanged) then A3,D3,A4,D4,A5,D5
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/6e684c28b587a70f0cfd5f5c6ff1ff3adb27a702/e2e/20bcb66c-dd1a-43b3-a248-76b294441b45.md", "", "", "20bcb66c-dd1a-43b3-a248-76b294441b45.md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/32900312ec1e6af4822ada052026ac7daaba561d/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/20bcb66c-dd1a-43b3-a248-76b294441b45.633f0c36d9a87b0eefdf682d49f1e52b4ef346d6.de-de.xlf", "", "", "20bcb66c-dd1a-43b3-a248-76b294441b45.633f0c36d9a87b0eefdf682d49f1e52b4ef346d6.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/9b8504e5aa1402fa50e4069711cfa2a1d738489d/e2e/20bcb66c-dd1a-43b3-a248-76b294441b45.md", "", "", "20bcb66c-dd1a-43b3-a248-76b294441b45.md")
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4a3e68fcbdc0c0002abf20f8990c30c3b09e6d57/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/20bcb66c-dd1a-43b3-a248-76b294441b45.633f0c36d9a87b0eefdf682d49f1e52b4ef346d6.de-de.xlf", "", "", "20bcb66c-dd1a-43b3-a248-76b294441b45.633f0c36d9a87b0eefdf682d49f1e52b4ef346d6.de-de.xlf")

$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/10bb3719a9cc4ae186d9c76b788262a8000000000/e2e/10bb3719-a9cc-4ae1-86d9-c76b788262a8.md", "", "", "10bb3719-a9cc-4ae1-86d9-c76b788262a8.md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b35df4560331639d247d8c3db9bcf659d9491bc4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/10bb3719-a9cc-4ae1-86d9-c76b788262a8.b35df4560331639d247d8c3db9bcf659d9491bc4.de-de.xlf", "", "", "10bb3719-a9cc-4ae1-86d9-c76b788262a8.b35df4560331639d247d8c3db9bcf659d9491bc4.de-de.xlf")

$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/66575674fd334ba19e478dae58941a75000000000/e2e/66575674-fd33-4ba1-9e47-8dae58941a75.md", "", "", "66575674-fd33-4ba1-9e47-8dae58941a75.md")
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2d11a3a5a73b03f5b090dc2164ab42aedd19db64/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/66575674-fd33-4ba1-9e47-8dae58941a75.2d11a3a5a73b03f5b090dc2164ab42aedd19db64.de-de.xlf", "", "", "66575674-fd33-4ba1-9e47-8dae58941a75.2d11a3a5a73b03f5b090dc2164ab42aedd19db64.de-de.xlf")

$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/4a68d5de10ad6c3300541df7be93d85829774398/e2e/b68a2137-e2a6-439b-88bb-be13c4dbe045.md", "", "", "b68a2137-e2a6-439b-88bb-be13c4dbe045.md")
$ws3.Hyperlinks.Add($ws3.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/52aaa75309f84d13edc3e411a90d4758a4ff3139/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b68a2137-e2a6-439b-88bb-be13c4dbe045.cafe1cadf14af19aedb092ec9298155126eb95ed.de-de.xlf", "", "", "b68a2137-e2a6-439b-88bb-be13c4dbe045.cafe1cadf14af19aedb092ec9298155126eb95ed.de-de.xlf")
